$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 397, shifting rows 397-502 down to 398-503.
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row 397 with its data.
$ws.Range("A397").Value = 10
$ws.Range("B397").Value = "Vega Modelo de Temuco"
$ws.Range("C397").Value = "La Araucanía"
$ws.Range("D397").Value = 44782
$ws.Range("E397").Value = 9
$ws.Range("F397").Value = "Fruta"
$ws.Range("G397").Value = 100101
$ws.Range("H397").Value = "Berries"
$ws.Range("I397").Value = 100101007
$ws.Range("J397").Value = "Kiwi"
$ws.Range("K397").Value = "Hayward"
$ws.Range("L397").Value = "Primera"
$ws.Range("M397").Value = 125
$ws.Range("N397").Value = 12000
$ws.Range("O397").Value = 12000
$ws.Range("P397").Value = 12000
$ws.Range("Q397").Value = "`$/bandeja 18 kilos"
$ws.Range("R397").Value = "Región de O'Higgins"
$ws.Range("S397").Value = 667
$ws.Range("T397").Value = 18
